# FIX: some fixes to the lexical scanner FSM.
#
# The "Transitions" sheet is a finite-state-machine transition table:
# row 1 / col A hold state & character-class labels, and each interior
# cell is the resulting state/token reached from the row's current state
# on the column's character class (or vice-versa, the layout is
# symmetric in this sheet). Column T is the "S_NUMBER" state (see T1).
# This edit fills in the previously-blank S_NUMBER column transitions,
# and the two new B4/B5 (S_START col) transitions for digit/zero.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transitions")

# --- S_START column (B): digit / zero now both enter S_NUMBER ---
$ws.Range("B4").Value = "S_NUMBER"
$ws.Range("B5").Value = "S_NUMBER"

# --- S_NUMBER column (T): terminal-token cells (bold, black) ---
# These delimiters close out a plain integer literal -> T_INTEGER,
# except C_PERCENT which closes out a percent literal -> T_PERCENT.
$intCells = "T2","T3","T6","T7","T8","T9","T10","T11","T12","T21","T25"
foreach ($addr in $intCells) {
    $cell = $ws.Range($addr)
    $cell.Value = "T_INTEGER"
    $cell.Font.Bold = $true
    $cell.Font.ColorIndex = -4105
}
$pctCell = $ws.Range("T23")
$pctCell.Value = "T_PERCENT"
$pctCell.Font.Bold = $true
$pctCell.Font.ColorIndex = -4105

# --- S_NUMBER column (T): continuation into another lexer state ---
# (normal, non-bold styling, matching the rest of the S_* cells)
$ws.Range("T4").Value  = "S_NUMBER"
$ws.Range("T5").Value  = "S_NUMBER"
$ws.Range("T13").Value = "S_SHARP"
$ws.Range("T14").Value = "S_NUMBER"
$ws.Range("T15").Value = "S_TIME_1ST"
$ws.Range("T16").Value = "S_PAIR_1ST"
$ws.Range("T17").Value = "S_DECIMAL"
$ws.Range("T19").Value = "S_DATE"
$ws.Range("T26").Value = "S_EMAIL"
$ws.Range("T27").Value = "S_DECIMAL"

# --- S_NUMBER column (T): illegal-character cells -> T_ERROR (bold red) ---
$errCells = "T18","T20","T22","T28","T29","T30","T32","T33","T35","T36"
foreach ($addr in $errCells) {
    $cell = $ws.Range($addr)
    $cell.Value = "T_ERROR"
    $cell.Font.Bold = $true
    $cell.Font.Color = 192
}

# T24 (C_COMMA x S_NUMBER) is also T_ERROR, but flagged with a yellow
# highlight fill to call out the known "1,000" decimal-comma mistake.
$flag = $ws.Range("T24")
$flag.Value = "T_ERROR"
$flag.Font.Bold = $true
$flag.Font.Color = 192
$flag.Interior.Color = 65535

# --- cosmetic: selection moved one cell right, now portrait page setup ---
[void]$ws.Range("T36").Select()
$ws.PageSetup.Orientation = 1
